$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers) ---
# Old P1 "Observación " / Q1 "Clave" shift right to S1 / T1 (same header style).
$oldP1 = $ws.Range("P1").Text
$oldQ1 = $ws.Range("Q1").Text

$ws.Range("P1").Copy() | Out-Null
$ws.Range("S1").PasteSpecial(-4122) | Out-Null
$ws.Range("Q1").Copy() | Out-Null
$ws.Range("T1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("S1").Value = $oldP1
$ws.Range("T1").Value = $oldQ1

# New P1, Q1, R1 become "Vacio" (keep the existing header style already on P1/Q1).
$ws.Range("P1").Value = "Vacio"
$ws.Range("Q1").Value = "Vacio"
$ws.Range("P1").Copy() | Out-Null
$ws.Range("R1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("R1").Value = "Vacio"

# --- Row 2 (data) ---
# Old P2 "Aprobado" / Q2 "prueba" shift right to S2 / T2, carrying their
# current (Text-formatted) style along.
$oldP2 = $ws.Range("P2").Text
$oldQ2 = $ws.Range("Q2").Text

$ws.Range("P2").Copy() | Out-Null
$ws.Range("S2").PasteSpecial(-4122) | Out-Null
$ws.Range("Q2").Copy() | Out-Null
$ws.Range("T2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("S2").Value = $oldP2
$ws.Range("T2").Value = $oldQ2

# New numeric flag columns P2, Q2, R2 = 1 (plain General number format,
# matching the row's own default style). Reset P2/Q2's leftover Text
# formatting by pasting the (General) format that R2 already carries
# before putting the new numbers in.
$ws.Range("R2").Copy() | Out-Null
$ws.Range("P2").PasteSpecial(-4122) | Out-Null
$ws.Range("Q2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 1

# U2 result code changes value, must stay text even though it looks numeric.
$ws.Range("U2").NumberFormat = "@"
$ws.Range("U2").Value = "4978844"

# New V2 cell with validation message.
$ws.Range("V2").Value = "Se han encontrado errores en la Validacion de la Propuesta"
